$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-04-25 Friday"

# Update each table cell by exact row/column position, assigning the
# cell Range.Text directly (content-addressed Find/Replace can cross-wire
# cells whose old/new text collide elsewhere in the table).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "40÷8=5, 0"
$t.Cell(1, 2).Range.Text = "97÷9=10, 7"
$t.Cell(1, 3).Range.Text = "14÷4=3, 2"
$t.Cell(1, 4).Range.Text = "27÷5=5, 2"
$t.Cell(1, 5).Range.Text = "62÷2=31, 0"

$t.Cell(5, 1).Range.Text = "26÷5=5, 1"
$t.Cell(5, 2).Range.Text = "96÷2=48, 0"
$t.Cell(5, 3).Range.Text = "85÷9=9, 4"
$t.Cell(5, 4).Range.Text = "50÷9=5, 5"
$t.Cell(5, 5).Range.Text = "26÷4=6, 2"

$t.Cell(9, 1).Range.Text = "41÷7=5, 6"
$t.Cell(9, 2).Range.Text = "11÷8=1, 3"
$t.Cell(9, 3).Range.Text = "71÷2=35, 1"
$t.Cell(9, 4).Range.Text = "36÷2=18, 0"
$t.Cell(9, 5).Range.Text = "69÷5=13, 4"

$t.Cell(13, 1).Range.Text = "90÷8=11, 2"
$t.Cell(13, 2).Range.Text = "34÷6=5, 4"
$t.Cell(13, 3).Range.Text = "92÷9=10, 2"
$t.Cell(13, 4).Range.Text = "58÷9=6, 4"
$t.Cell(13, 5).Range.Text = "39÷2=19, 1"

$t.Cell(17, 1).Range.Text = "66÷6=11, 0"
$t.Cell(17, 2).Range.Text = "20÷3=6, 2"
$t.Cell(17, 3).Range.Text = "24÷3=8, 0"
$t.Cell(17, 4).Range.Text = "96÷4=24, 0"
$t.Cell(17, 5).Range.Text = "49÷8=6, 1"
